$wb = $excel.ActiveWorkbook

# Fix typo in the "genotype" sheet's repeated label: b9d2_unkown -> b9d2_unknown
$genotype = $wb.Worksheets.Item("genotype")
$genotype.Range("B2:L9").Value = "b9d2_unknown"

# Make "temperature" tab no longer the selected/active tab, and make
# "genotype" the active tab instead, with cell B2 selected on it.
$temperature = $wb.Worksheets.Item("temperature")
$temperature.Range("A1:M9").Select()

$genotype.Activate()
$genotype.Range("B2").Select()
